# Apply the 2022-06-15 data update to the "Fonds de solidarite" VOLET1 dataset.
# Only the "nombre_aides" (column C) and "montant_total" (column E) values
# change for the rows listed below; "nombre_entreprises" (column D) is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 8;   C = 181370;  E = 653248410 },
    @{ Row = 10;  C = 278210;  E = 1752513740 },
    @{ Row = 19;  C = 108926;  E = 344804726 },
    @{ Row = 97;  C = 98510;   E = 307118902 },
    @{ Row = 152; C = 126054;  E = 716092217 },
    @{ Row = 164; C = 50586;   E = 168946548 },
    @{ Row = 168; C = 285122;  E = 1213790875 },
    @{ Row = 169; C = 562673;  E = 1286266968 },
    @{ Row = 170; C = 367577;  E = 2848223671 },
    @{ Row = 171; C = 115230;  E = 449025067 },
    @{ Row = 174; C = 357382;  E = 1020222665 },
    @{ Row = 175; C = 125700;  E = 815905290 },
    @{ Row = 179; C = 235806;  E = 813703438 },
    @{ Row = 180; C = 141532;  E = 341248758 },
    @{ Row = 279; C = 28968;   E = 57089102 },
    @{ Row = 293; C = 61673;   E = 194903805 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
